# Apply the "cryptos list" price/volume update described by the commit.
# Numeric-looking cells in column D must stay plain text (they are stored as
# inline strings in the source), so we prefix them with a leading apostrophe
# the same way a user typing into Excel would, to stop them turning into real
# numbers (which would silently drop formatting such as trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.330.67'   # was '51.287.34'
$ws.Range("E2").Value = '  +2.51%  '   # was '  +2.50%  '

$ws.Range("D3").Value = '2.745.51'   # was '2.745.70'
$ws.Range("E3").Value = '  +2.76%  '   # was '  +2.97%  '

$ws.Range("E4").Value = '  +0.02%  '   # was '  +0.12%  '

$ws.Range("D5").Value = "'" + '115.29'   # was '115.07'
$ws.Range("E5").Value = '  +0.96%  '   # was '  +0.90%  '

$ws.Range("D6").Value = "'" + '330.99'   # was '330.98'
$ws.Range("E6").Value = '  +1.24%  '   # was '  +1.32%  '

$ws.Range("E7").Value = '  +0.51%  '   # was '  +0.48%  '

$ws.Range("D8").Value = "'" + '0.999'   # was '1.00'
$ws.Range("E8").Value = '  -0.05%  '   # was '  +0.03%  '

$ws.Range("D9").Value = "'" + '0.571'   # was '0.568'
$ws.Range("E9").Value = '  +2.58%  '   # was '  +1.93%  '

$ws.Range("D10").Value = "'" + '41.43'   # was '41.41'
$ws.Range("E10").Value = '  +0.76%  '   # was '  +0.80%  '

$ws.Range("D11").Value = "'" + '20.29'   # was '20.27'
$ws.Range("E11").Value = '  +0.89%  '   # was '  +0.82%  '

$ws.Range("D12").Value = "'" + '0.0829'   # was '0.0828'
$ws.Range("E12").Value = '  +0.42%  '   # was '  +0.34%  '

$ws.Range("E13").Value = '  +2.70%  '   # was '  +2.84%  '

$ws.Range("E14").Value = '  +4.09%  '   # was '  +4.25%  '

$ws.Range("D15").Value = '3.175.75'   # was '3.183.43'
$ws.Range("E15").Value = '  +2.94%  '   # was '  +3.02%  '

$ws.Range("D16").Value = '2.743.60'   # was '2.740.11'
$ws.Range("E16").Value = '  +2.69%  '   # was '  +2.31%  '

$ws.Range("E17").Value = '  +0.85%  '   # was '  +0.87%  '

$ws.Range("D18").Value = '51.237.59'   # was '51.343.90'
$ws.Range("E18").Value = '  +2.48%  '   # was '  +2.74%  '

$ws.Range("D19").Value = "'" + '13.60'   # was '13.56'
$ws.Range("E19").Value = '  +2.64%  '   # was '  +2.26%  '

$ws.Range("D20").Value = "'" + '3.02'   # was '3.00'
$ws.Range("E20").Value = '  +4.52%  '   # was '  +3.81%  '

$ws.Range("D21").Value = "'" + '6.85'   # was '6.83'
$ws.Range("E21").Value = '  +0.85%  '   # was '  +0.60%  '

$ws.Range("D22").Value = '0.0₃0964'   # was '0.0₃0961'
$ws.Range("E22").Value = '  +0.14%  '   # was '  -0.03%  '

$ws.Range("D23").Value = "'" + '283.57'   # was '285.11'
$ws.Range("E23").Value = '  +2.03%  '   # was '  +2.67%  '

$ws.Range("D24").Value = "'" + '70.34'   # was '70.31'
$ws.Range("E24").Value = '  -3.14%  '   # was '  -3.18%  '

$ws.Range("E25").Value = '  +0.13%  '   # was '  +0.18%  '

$ws.Range("E26").Value = '  -0.32%  '   # was '  -0.29%  '

$ws.Range("D27").Value = "'" + '0.999'   # was '1.00'
$ws.Range("E27").Value = '  -0.03%  '   # was '  +0.07%  '

$ws.Range("E28").Value = '  +1.86%  '   # was '  +1.94%  '

$ws.Range("D29").Value = "'" + '2.23'   # was '2.22'
$ws.Range("E29").Value = '  -0.34%  '   # was '  -0.45%  '

$ws.Range("E30").Value = '  -1.90%  '   # was '  -1.76%  '

$ws.Range("D31").Value = "'" + '35.52'   # was '35.36'
$ws.Range("E31").Value = '  -2.74%  '   # was '  -3.32%  '

$ws.Range("D32").Value = "'" + '50.19'   # was '50.15'
$ws.Range("E32").Value = '  -0.30%  '   # was '  -0.29%  '

$ws.Range("D33").Value = "'" + '5.65'   # was '5.61'
$ws.Range("E33").Value = '  +2.67%  '   # was '  +2.14%  '

$ws.Range("D34").Value = "'" + '0.0825'   # was '0.0823'
$ws.Range("E34").Value = '  +0.79%  '   # was '  +0.53%  '

$ws.Range("D35").Value = "'" + '19.41'   # was '19.38'
$ws.Range("E35").Value = '  -1.66%  '   # was '  -1.75%  '

$ws.Range("E36").Value = '  -0.10%  '   # was '  -0.13%  '

$ws.Range("E37").Value = '  +1.08%  '   # was '  +1.14%  '

$ws.Range("E38").Value = '  -1.90%  '   # was '  -1.66%  '

$ws.Range("D39").Value = "'" + '3.22'   # was '3.21'
$ws.Range("E39").Value = '  +2.75%  '   # was '  +2.51%  '

$ws.Range("D40").Value = "'" + '129.24'   # was '129.56'
$ws.Range("E40").Value = '  +3.17%  '   # was '  +3.45%  '

$ws.Range("D41").Value = "'" + '23.67'   # was '23.69'
$ws.Range("E41").Value = '  +3.67%  '   # was '  +4.00%  '

$ws.Range("D42").Value = "'" + '0.0354'   # was '0.0352'
$ws.Range("E42").Value = '  +10.79%  '   # was '  +10.47%  '

$ws.Range("B43").Value = 'WEMIXToken'   # was 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'   # was 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = "'" + '2.29'   # was '0.113'
$ws.Range("E43").Value = '  +3.48%  '   # was '  +0.21%  '

$ws.Range("B44").Value = 'Stellar'   # was 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'   # was 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Value = "'" + '0.113'   # was '2.29'
$ws.Range("E44").Value = '  +0.26%  '   # was '  +3.01%  '

$ws.Range("D45").Value = "'" + '3.42'   # was '3.40'
$ws.Range("E45").Value = '  +2.91%  '   # was '  +2.66%  '

$ws.Range("D46").Value = '2.106.84'   # was '2.108.07'
$ws.Range("E46").Value = '  -0.38%  '   # was '  -0.35%  '

$ws.Range("D47").Value = "'" + '2.23'   # was '2.21'
$ws.Range("E47").Value = '  +10.33%  '   # was '  +9.61%  '

$ws.Range("E48").Value = '  -1.11%  '   # was '  -1.17%  '

$ws.Range("E49").Value = '  +2.12%  '   # was '  +2.13%  '

$ws.Range("D50").Value = "'" + '9.07'   # was '9.05'
$ws.Range("E50").Value = '  -0.56%  '   # was '  -0.77%  '

$ws.Range("D51").Value = "'" + '59.94'   # was '59.89'
$ws.Range("E51").Value = '  +0.03%  '   # was '  +0.04%  '
